$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.410.70"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "2.007.07"
$ws.Range("E3").Value = "  +5.56%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'244.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("E6").Value = "  -4.76%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'44.43"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("D9").Value = "'62.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.91%  "
$ws.Range("D10").Value = "'0.364"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("D11").Value = "'0.0712"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.14%  "
$ws.Range("D12").Value = "'0.0978"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").Value = "'14.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "2.296.52"
$ws.Range("E14").Value = "  +5.50%  "
$ws.Range("D15").Value = "'0.806"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "2.008.27"
$ws.Range("E16").Value = "  +6.23%  "
$ws.Range("E17").Value = "  -3.18%  "
$ws.Range("D18").Value = "36.343.70"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").Value = "'71.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.80%  "
$ws.Range("D20").Value = "0.0₃0813"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("D21").Value = "'12.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").Value = "'236.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.14%  "
$ws.Range("D23").Value = "'4.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.81%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'2.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.13%  "
$ws.Range("D26").Value = "'164.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("D27").Value = "'8.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "'19.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.29%  "
$ws.Range("D29").Value = "'1.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -11.45%  "
$ws.Range("E30").Value = "  -5.97%  "
$ws.Range("D31").Value = "'22.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +62.41%  "
$ws.Range("D32").Value = "'4.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("E33").Value = "  -3.74%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").Value = "'3.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.21%  "
$ws.Range("D37").Value = "'0.0820"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.67%  "
$ws.Range("E38").Value = "  +7.55%  "
$ws.Range("D39").Value = "'0.854"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "'1.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.76%  "
$ws.Range("E41").Value = "  -5.13%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'95.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("D44").Value = "'2.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +15.77%  "
$ws.Range("E45").Value = "  -7.32%  "
$ws.Range("D46").Value = "1.308.73"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("D47").Value = "'0.0815"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").Value = "2.189.37"
$ws.Range("E49").Value = "  +5.36%  "
$ws.Range("E50").Value = "  -8.32%  "
$ws.Range("E51").Value = "  +13.80%  "
